# Add new columns I (I0) and J (IF) to the worksheet, matching the
# formatting of the existing header/data columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new headers "I0" and "IF" with the same style as
# the existing header cells (bold, centered, bordered).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2-11 for columns I and J.
$values = @{
    2  = @(7, 9)
    3  = @(8, 9)
    4  = @(6, 8)
    5  = @(1, 5)
    6  = @(1, 5)
    7  = @(1, 4)
    8  = @(5, 9)
    9  = @(1, 7)
    10 = @(1, 2)
    11 = @(1, 5)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
